# Update "want-to-go" counts (column F) that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 47
$ws1.Range("F25").Value = 243
$ws1.Range("F28").Value = 6
$ws1.Range("F29").Value = 840
$ws1.Range("F30").Value = 77
$ws1.Range("F32").Value = 46

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F21").Value = 45
$ws2.Range("F31").Value = 10
$ws2.Range("F45").Value = 756

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F28").Value = 45
$ws4.Range("F32").Value = 243
$ws4.Range("F38").Value = 840
$ws4.Range("F40").Value = 77
$ws4.Range("F42").Value = 10
$ws4.Range("F52").Value = 756
